$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

$data = @{
    2  = @(0.2066981850894128, 2.270155390337222, 13.1891407897574, 3.631685667807362, 3.661877248494126, 51)
    3  = @(-0.075228059080023, 2.191059243838117, 12.66358371323103, 3.558592940086155, 3.59391839656686, 50)
    4  = @(0.207102134124417, 2.274418628387441, 13.76964981448061, 3.710747878053778, 3.743358472864432, 49)
    5  = @(0.03062386419637815, 2.209507036589459, 12.80187705469571, 3.577971080751731, 3.615701818853944, 48)
    6  = @(0.2145225631768136, 2.299235277238912, 13.98594729776688, 3.739779043976646, 3.77398585940008, 47)
    7  = @(0.1255984722002377, 2.192001147543182, 12.97327225105472, 3.601842896498225, 3.639428754450474, 46)
    8  = @(0.1386258507452061, 2.134562903283757, 12.94582025312243, 3.598030051725864, 3.635985363188866, 45)
    9  = @(-0.01227215085340897, 1.819249639236372, 11.02222735225344, 3.319973998731532, 3.358333532090894, 44)
    10 = @(-0.0200094225433071, 1.878919759582091, 11.5868556758706, 3.403947073012534, 3.44417236534218, 43)
    11 = @(0.06995512137203723, 1.589308054450445, 10.76865959017209, 3.28156358923183, 3.320586807489363, 42)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
